$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Simran Bala", "18", "Female", "05-01-2025 22:09:18", "ADHD", "attention deficiet stuff", "No", "Yes yes yes yes yes"),
    @("Subhra Bala", "54", "Female", "05-01-2025 22:32:23", "Eye pain, Heart Ache", "Uterus problem", "No", "Yes"),
    @("Guddu", "19", "Male", "06-01-2025 22:54:52", "No", "None", "Yes yes", "No NO"),
    @("Souradip Banerjee", "35", "Male", "08-01-2025 18:12:44", "Gastritis with sugar problem", "Heart infection", "Covid-20 Checkup", "eat healthy and workout")
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($c = 0; $c -lt 8; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$i][$c]
        $cell.Style = "Normal"
    }
}

for ($c = 1; $c -le 8; $c++) {
    $ws.Columns.Item($c).ClearFormats()
}

$ws.Range("L8").Select()
